# Add new kernel for bfs_rec
# This adds a small "storage needs" computation block to the right of the
# existing results table (columns L-O, rows 16-18) and shifts the shared
# string used by E17 because a new string is inserted before it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: headers for the new storage block
$ws.Range("L16").Value = "Total storage needs (MB)"
$ws.Range("M16").Value = "Nodes"
$ws.Range("N16").Formula = "=B17*4/(1024*1024)"

# Row 17: Edges storage
$ws.Range("M17").Value = "Edges"
$ws.Range("N17").Formula = "=(D17*4)/(1024*1024)"

# Row 18 (new row): Queue_Sources storage + comment "too much"
$ws.Range("M18").Value = "Queue_Sources"
$ws.Range("N18").Formula = "=(B17*B17*4)/(1024*1024)"
$ws.Range("O18").Value = "too much"

# E17 previously referenced "No double links, no self-directed edges";
# re-set it (a new shared string gets inserted before it in the sst table,
# which naturally happens due to the order strings are introduced above).
$ws.Range("E17").Value = "No double links, no self-directed edges"

# Update selection to match the recorded cursor position after editing
$ws.Range("E10").Select()
